$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 356 (old rows 356-441 shift down to 358-443)
$ws.Rows.Item(356).Insert()
$ws.Rows.Item(356).Insert()

# New row 356 data
$ws.Cells.Item(356, 1).Value2 = 3
$ws.Cells.Item(356, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(356, 3).Value2 = "Coquimbo"
$ws.Cells.Item(356, 4).Value2 = 44511
$ws.Cells.Item(356, 5).Value2 = 5
$ws.Cells.Item(356, 6).Value2 = 100112023
$ws.Cells.Item(356, 7).Value2 = "Brócoli"
$ws.Cells.Item(356, 8).Value2 = "Sin especificar"
$ws.Cells.Item(356, 9).Value2 = "Primera"
$ws.Cells.Item(356, 10).Value2 = 2588
$ws.Cells.Item(356, 11).Value2 = 600
$ws.Cells.Item(356, 12).Value2 = 650
$ws.Cells.Item(356, 13).Value2 = 625
$ws.Cells.Item(356, 14).Value2 = "$/unidad"
$ws.Cells.Item(356, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(356, 16).Value2 = 625
$ws.Cells.Item(356, 17).Value2 = 1
$ws.Cells.Item(356, 18).Value2 = "Hortaliza"

# New row 357 data
$ws.Cells.Item(357, 1).Value2 = 3
$ws.Cells.Item(357, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(357, 3).Value2 = "Coquimbo"
$ws.Cells.Item(357, 4).Value2 = 44511
$ws.Cells.Item(357, 5).Value2 = 5
$ws.Cells.Item(357, 6).Value2 = 100112023
$ws.Cells.Item(357, 7).Value2 = "Brócoli"
$ws.Cells.Item(357, 8).Value2 = "Sin especificar"
$ws.Cells.Item(357, 9).Value2 = "Segunda"
$ws.Cells.Item(357, 10).Value2 = 1380
$ws.Cells.Item(357, 11).Value2 = 500
$ws.Cells.Item(357, 12).Value2 = 500
$ws.Cells.Item(357, 13).Value2 = 500
$ws.Cells.Item(357, 14).Value2 = "$/unidad"
$ws.Cells.Item(357, 15).Value2 = "Provincia de Quillota"
$ws.Cells.Item(357, 16).Value2 = 500
$ws.Cells.Item(357, 17).Value2 = 1
$ws.Cells.Item(357, 18).Value2 = "Hortaliza"
